$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 91.5
$ws.Cells.Item(2, 9).Value = 91.10714
$ws.Cells.Item(2, 10).Value = 97
$ws.Cells.Item(2, 11).Value = 91.10714
$ws.Cells.Item(2, 12).Value = 97
$ws.Cells.Item(2, 13).Value = 21.89286
$ws.Cells.Item(2, 14).Value = -323
$ws.Cells.Item(21, 8).Value = 20750
$ws.Cells.Item(21, 9).Value = 14333.333
$ws.Cells.Item(21, 10).Value = 40000
$ws.Cells.Item(21, 11).Value = 14333.333
$ws.Cells.Item(21, 12).Value = 40000
$ws.Cells.Item(21, 13).Value = -13865.333
$ws.Cells.Item(21, 14).Value = -40936
$ws.Cells.Item(23, 8).Value = 20750
$ws.Cells.Item(23, 9).Value = 14333.333
$ws.Cells.Item(23, 10).Value = 40000
$ws.Cells.Item(23, 11).Value = 14333.333
$ws.Cells.Item(23, 12).Value = 40000
$ws.Cells.Item(23, 13).Value = -14099.333
$ws.Cells.Item(23, 14).Value = -40468
$ws.Cells.Item(29, 8).Value = 2975.3
$ws.Cells.Item(29, 9).Value = 50
$ws.Cells.Item(29, 10).Value = 3300.3333
$ws.Cells.Item(29, 11).Value = 150
$ws.Cells.Item(29, 12).Value = 9900.999899999999
$ws.Cells.Item(29, 13).Value = 131
$ws.Cells.Item(29, 14).Value = -10462.9999
$ws.Cells.Item(32, 8).Value = 757.8333
$ws.Cells.Item(32, 9).Value = 928.2857
$ws.Cells.Item(32, 10).Value = 519.2
$ws.Cells.Item(32, 11).Value = 928.2857
$ws.Cells.Item(32, 12).Value = 519.2
$ws.Cells.Item(32, 13).Value = -602.2857
$ws.Cells.Item(32, 14).Value = -1171.2
$ws.Cells.Item(38, 8).Value = 394.85715
$ws.Cells.Item(38, 9).Value = 152.8
$ws.Cells.Item(58, 8).Value = 1733.3448
$ws.Cells.Item(58, 9).Value = 353.4
$ws.Cells.Item(58, 10).Value = 2020.8334
$ws.Cells.Item(58, 11).Value = 1060.2
$ws.Cells.Item(58, 12).Value = 6062.5002
$ws.Cells.Item(58, 13).Value = -910.1999999999998
$ws.Cells.Item(58, 14).Value = -6362.5002
$ws.Cells.Item(87, 8).Value = 35600
$ws.Cells.Item(87, 9).Value = 40000
$ws.Cells.Item(87, 10).Value = 34500
$ws.Cells.Item(87, 11).Value = 40000
$ws.Cells.Item(87, 12).Value = 34500
$ws.Cells.Item(87, 13).Value = -38752
$ws.Cells.Item(87, 14).Value = -36996
$ws.Cells.Item(90, 8).Value = 35600
$ws.Cells.Item(90, 9).Value = 40000
$ws.Cells.Item(90, 10).Value = 34500
$ws.Cells.Item(90, 11).Value = 120000
$ws.Cells.Item(90, 12).Value = 103500
$ws.Cells.Item(90, 13).Value = -113760
$ws.Cells.Item(90, 14).Value = -115980
$ws.Cells.Item(107, 8).Value = 748.7692
$ws.Cells.Item(107, 9).Value = 236.3077
$ws.Cells.Item(107, 10).Value = 1261.2307
$ws.Cells.Item(107, 11).Value = 236.3077
$ws.Cells.Item(107, 12).Value = 1261.2307
$ws.Cells.Item(107, 13).Value = 1683.6923
$ws.Cells.Item(107, 14).Value = -5101.2307
$ws.Cells.Item(113, 8).Value = 3179.6
$ws.Cells.Item(113, 9).Value = 3250
$ws.Cells.Item(113, 10).Value = 3162
$ws.Cells.Item(113, 11).Value = 3250
$ws.Cells.Item(113, 12).Value = 3162
$ws.Cells.Item(113, 13).Value = 4
$ws.Cells.Item(113, 14).Value = -9670
$ws.Cells.Item(116, 8).Value = 9953
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 9953
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = ""
$ws.Cells.Item(116, 13).Value = 9953
$ws.Cells.Item(116, 14).Value = -16837
$ws.Cells.Item(127, 8).Value = 690.3333
$ws.Cells.Item(127, 9).Value = 671.2727
$ws.Cells.Item(127, 10).Value = 900
$ws.Cells.Item(127, 11).Value = 2013.8181
$ws.Cells.Item(127, 12).Value = 2700
$ws.Cells.Item(127, 13).Value = 2946.1819
$ws.Cells.Item(127, 14).Value = -12620
$ws.Cells.Item(137, 8).Value = 27028582
$ws.Cells.Item(137, 9).Value = 1243.4193
$ws.Cells.Item(137, 10).Value = 166669840
$ws.Cells.Item(137, 11).Value = 3730.2579
$ws.Cells.Item(137, 12).Value = 500009520
$ws.Cells.Item(137, 13).Value = -1180.2579
$ws.Cells.Item(137, 14).Value = -500014620
$ws.Cells.Item(138, 8).Value = 2033.6063
$ws.Cells.Item(138, 9).Value = 802.0345
$ws.Cells.Item(138, 10).Value = 2583.077
$ws.Cells.Item(138, 11).Value = 2406.1035
$ws.Cells.Item(138, 12).Value = 7749.231000000001
$ws.Cells.Item(138, 13).Value = 2733.8965
$ws.Cells.Item(138, 14).Value = -18029.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 38788.62
$ws.Cells.Item(32, 9).Value = 51602.582
$ws.Cells.Item(32, 10).Value = 26960.346
$ws.Cells.Item(32, 11).Value = 51602.582
$ws.Cells.Item(32, 12).Value = 26960.346
$ws.Cells.Item(32, 13).Value = -51315.582
$ws.Cells.Item(32, 14).Value = -27534.346
$ws.Cells.Item(61, 8).Value = 1722.15
$ws.Cells.Item(61, 9).Value = 1355.375
$ws.Cells.Item(61, 10).Value = 3189.25
$ws.Cells.Item(61, 11).Value = 1355.375
$ws.Cells.Item(61, 12).Value = 3189.25
$ws.Cells.Item(61, 13).Value = -1143.375
$ws.Cells.Item(61, 14).Value = -3613.25
$ws.Cells.Item(74, 8).Value = 1345.4
$ws.Cells.Item(74, 9).Value = 1296
$ws.Cells.Item(74, 10).Value = 1398.9166
$ws.Cells.Item(74, 11).Value = 1296
$ws.Cells.Item(74, 12).Value = 1398.9166
$ws.Cells.Item(74, 13).Value = -422
$ws.Cells.Item(74, 14).Value = -3146.9166
$ws.Cells.Item(77, 8).Value = 1345.4
$ws.Cells.Item(77, 9).Value = 1296
$ws.Cells.Item(77, 10).Value = 1398.9166
$ws.Cells.Item(77, 11).Value = 6480
$ws.Cells.Item(77, 12).Value = 6994.583000000001
$ws.Cells.Item(77, 13).Value = -2112
$ws.Cells.Item(77, 14).Value = -15730.583
$ws.Cells.Item(102, 8).Value = 1310
$ws.Cells.Item(102, 9).Value = 1233.3334
$ws.Cells.Item(102, 10).Value = 2000
$ws.Cells.Item(102, 11).Value = 1233.3334
$ws.Cells.Item(102, 12).Value = 2000
$ws.Cells.Item(102, 13).Value = 388.6666
$ws.Cells.Item(102, 14).Value = -5244
$ws.Cells.Item(132, 8).Value = 6219.2964
$ws.Cells.Item(132, 9).Value = 8044.758
$ws.Cells.Item(132, 10).Value = 3350.7144
$ws.Cells.Item(132, 11).Value = 24134.274
$ws.Cells.Item(132, 12).Value = 10052.1432
$ws.Cells.Item(132, 13).Value = -21604.274
$ws.Cells.Item(132, 14).Value = -15112.1432
$ws.Cells.Item(136, 8).Value = 1722.15
$ws.Cells.Item(136, 9).Value = 1355.375
$ws.Cells.Item(136, 10).Value = 3189.25
$ws.Cells.Item(136, 11).Value = 4066.125
$ws.Cells.Item(136, 12).Value = 9567.75
$ws.Cells.Item(136, 13).Value = -1516.125
$ws.Cells.Item(136, 14).Value = -14667.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 381.4
$ws.Cells.Item(80, 9).Value = 556
$ws.Cells.Item(80, 10).Value = 337.75
$ws.Cells.Item(80, 11).Value = 556
$ws.Cells.Item(80, 12).Value = 337.75
$ws.Cells.Item(80, 13).Value = 442
$ws.Cells.Item(80, 14).Value = -2333.75
$ws.Cells.Item(83, 8).Value = 381.4
$ws.Cells.Item(83, 9).Value = 556
$ws.Cells.Item(83, 10).Value = 337.75
$ws.Cells.Item(83, 11).Value = 2780
$ws.Cells.Item(83, 12).Value = 1688.75
$ws.Cells.Item(83, 13).Value = 2212
$ws.Cells.Item(83, 14).Value = -11672.75
$ws.Cells.Item(99, 8).Value = 935.1667
$ws.Cells.Item(99, 9).Value = 925
$ws.Cells.Item(99, 10).Value = 955.5
$ws.Cells.Item(99, 11).Value = 925
$ws.Cells.Item(99, 12).Value = 955.5
$ws.Cells.Item(99, 13).Value = 573
$ws.Cells.Item(99, 14).Value = -3951.5
$ws.Cells.Item(107, 8).Value = 5000
$ws.Cells.Item(107, 9).Value = 5000
$ws.Cells.Item(107, 10).Value = 5000
$ws.Cells.Item(107, 11).Value = 5000
$ws.Cells.Item(107, 12).Value = 5000
$ws.Cells.Item(107, 13).Value = -3080
$ws.Cells.Item(107, 14).Value = -8840

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 388.77777
$ws.Cells.Item(22, 9).Value = 374.875
$ws.Cells.Item(22, 10).Value = 500
$ws.Cells.Item(22, 11).Value = 374.875
$ws.Cells.Item(22, 12).Value = 500
$ws.Cells.Item(22, 13).Value = -24.875
$ws.Cells.Item(22, 14).Value = -1200
$ws.Cells.Item(132, 8).Value = 3206733.5
$ws.Cells.Item(132, 9).Value = 1263.72
$ws.Cells.Item(132, 10).Value = 8930787
$ws.Cells.Item(132, 11).Value = 3791.16
$ws.Cells.Item(132, 12).Value = 26792361
$ws.Cells.Item(132, 13).Value = -1261.16
$ws.Cells.Item(132, 14).Value = -26797421

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(101, 8).Value = 22833.334
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 22833.334
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 68500.00199999999
$ws.Cells.Item(101, 14).Value = -73368.00199999999
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 12).Value = ""
$ws.Cells.Item(102, 13).Value = ""
$ws.Cells.Item(102, 14).Value = 0
$ws.Cells.Item(103, 8).Value = 4107.2
$ws.Cells.Item(103, 9).Value = 800
$ws.Cells.Item(103, 10).Value = 4343.4287
$ws.Cells.Item(103, 11).Value = 2400
$ws.Cells.Item(103, 12).Value = 13030.2861
$ws.Cells.Item(103, 13).Value = -1521
$ws.Cells.Item(103, 14).Value = -14788.2861
$ws.Cells.Item(104, 8).Value = 4900.25
$ws.Cells.Item(104, 9).Value = 2026
$ws.Cells.Item(104, 10).Value = 5858.3335
$ws.Cells.Item(104, 11).Value = 6078
$ws.Cells.Item(104, 12).Value = 17575.0005
$ws.Cells.Item(104, 13).Value = -3457
$ws.Cells.Item(104, 14).Value = -22817.0005
$ws.Cells.Item(105, 8).Value = 10406.556
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 10406.556
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = ""
$ws.Cells.Item(105, 13).Value = 31219.668
$ws.Cells.Item(105, 14).Value = -36461.66800000001
$ws.Cells.Item(106, 8).Value = 6000
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 6000
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 18000
$ws.Cells.Item(106, 13).Value = 0
$ws.Cells.Item(131, 8).Value = 1441920
$ws.Cells.Item(131, 9).Value = 36043.332
$ws.Cells.Item(131, 10).Value = 1486788.4
$ws.Cells.Item(131, 11).Value = 108129.996
$ws.Cells.Item(131, 12).Value = 4460365.199999999
$ws.Cells.Item(131, 13).Value = -103089.996
$ws.Cells.Item(131, 14).Value = -4470445.199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = ""
$ws.Cells.Item(23, 14).Value = 0
$ws.Cells.Item(122, 8).Value = 24392590
$ws.Cells.Item(122, 9).Value = 40002220
$ws.Cells.Item(122, 10).Value = 2543.75
$ws.Cells.Item(122, 11).Value = 120006660
$ws.Cells.Item(122, 12).Value = 7631.25
$ws.Cells.Item(122, 13).Value = -120004210
$ws.Cells.Item(122, 14).Value = -12531.25
$ws.Cells.Item(132, 8).Value = 5868.8667
$ws.Cells.Item(132, 9).Value = 8606.625
$ws.Cells.Item(132, 10).Value = 2740
$ws.Cells.Item(132, 11).Value = 25819.875
$ws.Cells.Item(132, 12).Value = 8220
$ws.Cells.Item(132, 13).Value = -23289.875
$ws.Cells.Item(132, 14).Value = -13280

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2975.5557
$ws.Cells.Item(40, 9).Value = 2630
$ws.Cells.Item(40, 10).Value = 3666.6667
$ws.Cells.Item(40, 11).Value = 2630
$ws.Cells.Item(40, 12).Value = 3666.6667
$ws.Cells.Item(40, 13).Value = -2494
$ws.Cells.Item(40, 14).Value = -3938.6667
$ws.Cells.Item(61, 8).Value = 11495556
$ws.Cells.Item(61, 9).Value = 1121.5883
$ws.Cells.Item(61, 10).Value = 27779338
$ws.Cells.Item(61, 11).Value = 1121.5883
$ws.Cells.Item(61, 12).Value = 27779338
$ws.Cells.Item(61, 13).Value = -919.5882999999999
$ws.Cells.Item(61, 14).Value = -27779742
$ws.Cells.Item(113, 8).Value = 11495556
$ws.Cells.Item(113, 9).Value = 1121.5883
$ws.Cells.Item(113, 10).Value = 27779338
$ws.Cells.Item(113, 11).Value = 1121.5883
$ws.Cells.Item(113, 12).Value = 27779338
$ws.Cells.Item(113, 13).Value = 1048.4117
$ws.Cells.Item(113, 14).Value = -27783678
$ws.Cells.Item(122, 8).Value = 8100.5713
$ws.Cells.Item(122, 9).Value = 14934.667
$ws.Cells.Item(122, 10).Value = 2975
$ws.Cells.Item(122, 11).Value = 44804.001
$ws.Cells.Item(122, 12).Value = 8925
$ws.Cells.Item(122, 13).Value = -42354.001
$ws.Cells.Item(122, 14).Value = -13825

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 450
$ws.Cells.Item(100, 9).Value = 421.42856
$ws.Cells.Item(100, 10).Value = 550
$ws.Cells.Item(100, 11).Value = 842.85712
$ws.Cells.Item(100, 12).Value = 1100
$ws.Cells.Item(100, 13).Value = -301.85712
$ws.Cells.Item(100, 14).Value = -2182
$ws.Cells.Item(107, 8).Value = 15805047
$ws.Cells.Item(107, 9).Value = 6250429
$ws.Cells.Item(107, 10).Value = 37037532
$ws.Cells.Item(107, 11).Value = 18751287
$ws.Cells.Item(107, 12).Value = 111112596
$ws.Cells.Item(107, 13).Value = -18749367
$ws.Cells.Item(107, 14).Value = -111116436
$ws.Cells.Item(122, 8).Value = 2613.0417
$ws.Cells.Item(122, 9).Value = 2509.682
$ws.Cells.Item(122, 10).Value = 3750
$ws.Cells.Item(122, 11).Value = 7529.045999999999
$ws.Cells.Item(122, 12).Value = 11250
$ws.Cells.Item(122, 13).Value = -5079.045999999999
$ws.Cells.Item(122, 14).Value = -16150
$ws.Cells.Item(132, 8).Value = 2211
$ws.Cells.Item(132, 9).Value = 1935.5294
$ws.Cells.Item(132, 10).Value = 2679.3
$ws.Cells.Item(132, 11).Value = 5806.5882
$ws.Cells.Item(132, 12).Value = 8037.900000000001
$ws.Cells.Item(132, 13).Value = -3276.5882
$ws.Cells.Item(132, 14).Value = -13097.9
